$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.854.57"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").Value = "1.809.32"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'310.09"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'0.4484"
$ws.Range("E7").Value = "  +5.62%  "
$ws.Range("D8").Value = "'0.3670"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "'0.07281"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "'0.8522"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").Value = "'20.64"
$ws.Range("E11").Value = "  -2.32%  "
$ws.Range("D12").Value = "1.807.14"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "'6.600"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "'0.07080"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "'5.293"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "'91.45"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "'0.000008707"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "'14.84"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").Value = "26.867.47"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").Value = "'5.134"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'10.80"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "'1.978"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").Value = "'151.23"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").Value = "'2.235"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("D27").Value = "'18.39"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'5.195"
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("D29").Value = "'116.06"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "'0.08822"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").Value = "'1.172"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").Value = "'0.7466"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("D33").Value = "'2.927"
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").Value = "'4.427"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").Value = "'1.001"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").Value = "'1.084"
$ws.Range("E36").Value = "  -3.59%  "
$ws.Range("D37").Value = "'0.01957"
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").Value = "'0.05176"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").Value = "'0.5269"
$ws.Range("E39").Value = "  +3.50%  "
$ws.Range("D40").Value = "'2.864"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "'7.074"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("D42").Value = "'0.1686"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "'0.5188"
$ws.Range("E43").Value = "  +8.56%  "
$ws.Range("D44").Value = "'8.434"
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("D45").Value = "'10.52"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").Value = "'1.965"
$ws.Range("E46").Value = "  +6.67%  "
$ws.Range("D47").Value = "'105.17"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").Value = "'1.000"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").Value = "'1.651"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").Value = "'0.06312"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("D51").Value = "'0.9125"
$ws.Range("E51").Value = "  -0.82%  "
